# Update gh-pages to output generated at 456a3b4
# This updates the "想去人数" (F column) counts on the "展览" and "全部类型"
# worksheets to reflect newly scraped numbers.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @{
    3  = 12785
    4  = 26
    5  = 79
    6  = 59
    8  = 17
    9  = 8
    10 = 12703
    11 = 267
    12 = 17
    13 = 8291
    14 = 7624
    15 = 182
    18 = 118
    19 = 976
    22 = 372
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    4  = 12785
    5  = 26
    6  = 79
    7  = 59
    9  = 17
    10 = 8
    11 = 12703
    12 = 267
    13 = 17
    14 = 8293
    15 = 7624
    16 = 182
    19 = 118
    20 = 976
    24 = 372
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
